$wb = $excel.ActiveWorkbook

# The data sheet holding the boolean control lever
$ws = $wb.Worksheets.Item("BESHFoFRV")

# Set the control lever value from 0 (no foresight) to 1 (foresight)
$ws.Range("B2").Value = 1

# Make this sheet the active sheet and select B3 (matches the saved view state)
$ws.Activate()
$ws.Range("B3").Select()
